# CIERRE 16 FEB 22
# Applies the Feb-2022 credit-remission close: fills in the missing
# rows 18-30 (dates, store names, amounts) on the FEBRERO sheet,
# completes the F6/G6 payment-date/amount pair, and tweaks a couple of
# cosmetic view/format details on the ENERO sheet.

$wb = $excel.ActiveWorkbook

$wsEnero   = $wb.Worksheets.Item(4)   # REMISIONES   ENERO  2022
$wsFebrero = $wb.Worksheets.Item(5)   # REMISIONES FEBRERO   2022

# ---------------------------------------------------------------
# 1) ENERO sheet: H41 loses its highlight fill (style 71 -> 75-ish,
#    i.e. the "done" cell is no longer flagged in orange).
# ---------------------------------------------------------------
$wsEnero.Range("H41").Interior.Pattern = -4142   # xlPatternNone

# ---------------------------------------------------------------
# 2) FEBRERO sheet: row 6 payment-date / payment-amount filled in.
# ---------------------------------------------------------------
$wsFebrero.Range("F6").Value = 44601
$wsFebrero.Range("G6").Value = 45782

# ---------------------------------------------------------------
# 3) FEBRERO sheet: rows 18-30 get their remision date, client and
#    amount filled in (previously blank placeholders).
# ---------------------------------------------------------------
$rows = @(
    @{ Row = 18; Date = 44599; Client = "OBRADOR";          Amount = 6834  },
    @{ Row = 19; Date = 44599; Client = "COMERCIO CENTRAL "; Amount = 9261  },
    @{ Row = 20; Date = 44600; Client = "COMERCIO CENTRAL "; Amount = 23661 },
    @{ Row = 21; Date = 44601; Client = "ISRAEL LEDO";       Amount = 48706 },
    @{ Row = 22; Date = 44601; Client = "COMERCIO CENTRAL "; Amount = 16765 },
    @{ Row = 23; Date = 44601; Client = "OBRADOR";          Amount = 2040  },
    @{ Row = 24; Date = 44601; Client = "PROSUBCA";          Amount = 2655  },
    @{ Row = 25; Date = 44602; Client = "OBRADOR";          Amount = 6814  },
    @{ Row = 26; Date = 44602; Client = "COMERCIO CENTRAL "; Amount = 8547  },
    @{ Row = 27; Date = 44603; Client = "COMERCIO CENTRAL "; Amount = 35414 },
    @{ Row = 28; Date = 44604; Client = "COMERCIO CENTRAL "; Amount = 7035  },
    @{ Row = 29; Date = 44604; Client = "OBRADOR";          Amount = 75    },
    @{ Row = 30; Date = 44605; Client = "COMERCIO CENTRAL "; Amount = 6556  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $wsFebrero.Range("A$row").Value = $r.Date
    $wsFebrero.Range("D$row").Value = $r.Client
    $wsFebrero.Range("E$row").Value = $r.Amount
}

# ---------------------------------------------------------------
# 4) View/selection bookkeeping (matches the saved cursor position
#    in the workbook at close time). FEBRERO is reactivated last so
#    it remains the workbook's active tab, same as before the edit.
# ---------------------------------------------------------------
[void]$wsEnero.Activate()
[void]$wsEnero.Range("D42").Select()

[void]$wsFebrero.Activate()
[void]$wsFebrero.Range("E31").Select()
